# Auto-generated edit script: update cryptos list with new prices/volumes
# as captured in commit "Updated cryptos list on Wed Apr  5 20:53:31 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "1.001") must be
# forced to stay text, matching the workbook's original string cells.
$textCells = $excel.Union($ws.Range("D4"), $ws.Range("D5"), $ws.Range("D6"), $ws.Range("D7"), $ws.Range("D8"), $ws.Range("D9"), $ws.Range("D10"), $ws.Range("D11"), $ws.Range("D12"), $ws.Range("D13"), $ws.Range("D15"), $ws.Range("D16"), $ws.Range("D17"), $ws.Range("D18"), $ws.Range("D19"), $ws.Range("D20"), $ws.Range("D21"), $ws.Range("D22"), $ws.Range("D24"), $ws.Range("D25"), $ws.Range("D26"), $ws.Range("D29"), $ws.Range("D30"), $ws.Range("D31"), $ws.Range("D33"), $ws.Range("D34"), $ws.Range("D35"), $ws.Range("D36"), $ws.Range("D37"), $ws.Range("D38"), $ws.Range("D39"), $ws.Range("D40"), $ws.Range("D41"), $ws.Range("D42"), $ws.Range("D43"), $ws.Range("D44"), $ws.Range("D45"), $ws.Range("D46"), $ws.Range("D47"), $ws.Range("D49"), $ws.Range("D50"), $ws.Range("D51"))
foreach ($area in $textCells.Areas) { $area.NumberFormat = "@" }

# Apply all cell value updates
$ws.Range("D2").Value = "28.218.11"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.908.16"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "314.60"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "0.5066"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.3933"
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").Value = "0.09318"
$ws.Range("E9").Value = "  -5.59%  "
$ws.Range("D10").Value = "1.142"
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("D11").Value = "41.87"
$ws.Range("E11").Value = "  +2.37%  "
$ws.Range("D12").Value = "6.410"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "20.89"
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.910.83"
$ws.Range("E14").Value = "  +1.94%  "
$ws.Range("D15").Value = "7.306"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "0.00001124"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "92.64"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").Value = "0.06614"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "17.99"
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "6.205"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "28.280.08"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "11.44"
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("D25").Value = "2.320"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").Value = "2.600"
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("D27").Value = "2.133.32"
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("D29").Value = "158.05"
$ws.Range("D30").Value = "127.30"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").Value = "1.104"
$ws.Range("E31").Value = "  +3.14%  "
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("D33").Value = "5.648"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "3.617"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").Value = "9.674"
$ws.Range("E35").Value = "  +1.29%  "
$ws.Range("D36").Value = "0.06654"
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("D37").Value = "0.02418"
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").Value = "1.246"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "0.2187"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "1.269"
$ws.Range("E40").Value = "  +6.95%  "
$ws.Range("D41").Value = "0.6445"
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("D42").Value = "5.013"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").Value = "11.50"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").Value = "13.36"
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").Value = "0.6020"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "3.720"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("D49").Value = "2.019"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("D50").Value = "122.98"
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("D51").Value = "1.104"
$ws.Range("E51").Value = "  +3.14%  "

# Restore default (General) style on the forced-text cells so no stray
# formatting/style metadata is left behind
foreach ($area in $textCells.Areas) { $area.Style = "Normal" }

